$d = $word.ActiveDocument

# 1) "I like the organization because" -> new mission statement sentence.
$d.Content.Find.Execute(
    "I like the organization because", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "I want to be a part of company’s mission of ______ to operate, improve and safeguard our global trade ecosystem",
    2) | Out-Null

# 2) Delete the whole "Once a stable trade network can be disrupted by a geopolitican
#    events" bullet paragraph outright (its content is not reused anywhere).
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "geopolitican") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# 3) Extend "...I understand how more frustrating it is when unexpectancies occur."
#    with the new trailing clause about trade wars / COVID-19.
$d.Content.Find.Execute(
    " occur.", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " occur such as trade wars and COVID-19.",
    2) | Out-Null

# 4) Delete the trailing "I feel connected by the mission of streamlining the
#    process..." bullet paragraph outright.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "I feel connected by the mission of streamlining the process") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# 5) Clean up "...React to improve my " run split (drop the stray gramStart/gramEnd
#    proof-error split around "my").
$d.Content.Find.Execute(
    "After leaving the company, I taught myself React to improve my ", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "After leaving the company, I taught myself React to improve my ",
    2) | Out-Null
